$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two stray rows that were mistakenly added (ccrowa01/ccrowa02 +
# their IP addresses) by deleting B16:C17's contents and removing row 17
# entirely.
$ws.Range("B16:C16").ClearContents()
$ws.Rows("17").Delete()

# Restore the selection to B8 (matches the author's saved view).
$ws.Range("B8").Select()
